$d = $word.ActiveDocument

# Locate the end of the last paragraph in the document body (just before
# its paragraph mark) so the new paragraphs are appended after it and
# before the trailing sectPr.
$lastPara = $d.Paragraphs.Last
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

# Paragraph 1 (es-ES): "Minimos cuadrados aplicados a sistemas dinamicos..."
$p1 = '<w:p><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">Mínimos cuadrados aplicados a sistemas </w:t></w:r>' +
      '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">dinámicos se basa en tener datos basados en el tiempo y dar un modelo productivo para predecir </w:t></w:r>' +
      '</w:p>'

# Paragraph 2 (en-US): the "Yi=aY(i-2) + bY(i-1) + cYi(i-2)" formula,
# with proofErr spell-check markers bracketing the coined variable names.
$p2 = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Yi=</w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>aY</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>i-</w:t></w:r>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">2) + </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>bY</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(i-1)</w:t></w:r>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>+</w:t></w:r>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>cYi</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(i-2)</w:t></w:r>' +
      '</w:p>'

# Paragraph 3 (default formatting): "Calcular adecuadamente los coeficientes..."
$p3 = '<w:p>' +
      '<w:r><w:t xml:space="preserve">Calcular </w:t></w:r>' +
      '<w:r><w:t xml:space="preserve">adecuadamente </w:t></w:r>' +
      '<w:r><w:t xml:space="preserve">los coeficientes </w:t></w:r>' +
      '<w:r><w:t>para q</w:t></w:r>' +
      '<w:r><w:t>ue Y se mantenga igual pero qu</w:t></w:r>' +
      '<w:r><w:t>e</w:t></w:r>' +
      '<w:r><w:t xml:space="preserve"> nos de</w:t></w:r>' +
      '<w:r><w:t>l comportamiento en el tiempo</w:t></w:r>' +
      '</w:p>'

$body = $p1 + $p2 + $p3

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' + $body + '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($xml)
